# Apply trade #127 (closed at 2026-02-17 16:07:30) to the live trading
# results workbook and roll the new numbers up into the summary sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.46
$summary.Range("B4").Value = -0.55
$summary.Range("B5").Value = -0.09
$summary.Range("B6").Value = 127
$summary.Range("B8").Value = 62
$summary.Range("B9").Value = 37.01

# ---------------------------------------------------------------------
# Strategy Status sheet (row 4 = MarketMaking)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.45999999999999
$status.Range("D4").Value = 127
$status.Range("E4").Value = -0.55
$status.Range("F4").Value = -0.54
$status.Range("G4").Value = 37.01

# ---------------------------------------------------------------------
# All Trades + MarketMaking sheets: append the new trade as row 128.
# Copy row 127 first so formatting/types (e.g. the plain-text date and
# time strings) are preserved instead of being re-interpreted by Excel's
# autoformatting, then overwrite the cells that actually changed.
# ---------------------------------------------------------------------
$sheetNames = @("All Trades", "MarketMaking")
foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("A127:Q127").Copy($ws.Range("A128:Q128"))

    # Only touch the cells whose value actually differs from row 127 so
    # that the copied plain-text cells (date/time strings, labels) are
    # left alone and never re-interpreted by Excel's autoformatting.
    $ws.Cells.Item(128, 1).Value = 127
    $ws.Cells.Item(128, 3).Value = "16:07:30"
    $ws.Cells.Item(128, 7).Value = 0.19
    $ws.Cells.Item(128, 9).Value = -9.5238
    $ws.Cells.Item(128, 10).Value = -0.02
    $ws.Cells.Item(128, 11).Value = 99.45999999999999
}
